$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing A:E data to B:F
$ws.Range("A1").EntireColumn.Insert()

# Copy the header style (from B1, which used to be A1) to the new A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats

# Set header and row label values for the new ID column
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = "Hb 20"
$ws.Range("A3").Value = "S 17"
$ws.Range("A4").Value = "Hb 30"
$ws.Range("A5").Value = "KHb 25"
$ws.Range("A6").Value = "Hb 25"
$ws.Range("A7").Value = "KHb 30"
$ws.Range("A8").Value = "KHb 31"
$ws.Range("A9").Value = "KS 76"
